$p = $ppt.ActivePresentation

$oldDate = "2023-10-05"
$newDate = "10/01/2023"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master date placeholder
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# All slide layouts' date placeholders
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

# Move the "Oval 20" shape on slide 1 to its new position
# target OOXML offset: x=6895365 EMU, y=519915 EMU (EMU / 12700 = points)
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Oval 20") {
        $shp.Left = 6895365 / 12700.0
        $shp.Top = 519915 / 12700.0
    }
}
